$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit permutes the data rows (2, 4-9) of the sheet: each target row
# receives the full row content (columns A-AC that vary) that previously
# lived in a different source row. Row 3 is untouched by the edit.
#
# Mapping: target row -> source row (content to copy)
#   2 <- 7
#   4 <- 8
#   5 <- 9
#   6 <- 5
#   7 <- 2
#   8 <- 6
#   9 <- 4

# First snapshot the "before" values for every source row/column that we need,
# since writes to earlier rows must not clobber values still needed later.
$cols = @("A","B","D","E","F","G","H","Q","R","AC")

$snapshot = @{}
foreach ($r in 2,4,5,6,7,8,9) {
    $snapshot[$r] = @{}
    foreach ($col in $cols) {
        $snapshot[$r][$col] = $ws.Range("$col$r").Value2
    }
}

$mapping = @{ 2 = 7; 4 = 8; 5 = 9; 6 = 5; 7 = 2; 8 = 6; 9 = 4 }

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    foreach ($col in $cols) {
        $val = $snapshot[$sourceRow][$col]
        if ($col -eq "AC" -and ($null -eq $val)) {
            $ws.Range("AC$targetRow").Value2 = ""
        } else {
            $ws.Range("$col$targetRow").Value2 = $val
        }
    }
}
